$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80 (shifts existing rows 80-100 down to 81-101)
$ws.Rows.Item(80).Insert()

# Populate the new row 80 with the new weekly price record
$ws.Range("A80").Value = 10
$ws.Range("B80").Value = "Vega Modelo de Temuco"
$ws.Range("C80").Value = "La Araucanía"
$ws.Range("D80").Value = 44951
$ws.Range("E80").Value = 9
$ws.Range("F80").Value = 100112030
$ws.Range("G80").Value = "Poroto granado"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 55
$ws.Range("K80").Value = 45000
$ws.Range("L80").Value = 45000
$ws.Range("M80").Value = 45000
$ws.Range("N80").Value = "$/saco 25 kilos"
$ws.Range("O80").Value = "Región del Maule"
$ws.Range("P80").Value = 1800
$ws.Range("Q80").Value = 25
$ws.Range("R80").Value = "Hortaliza"
